$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2992
$ws.Range("C2").Value = 70.75

$ws.Range("B3").Value = 1017
$ws.Range("C3").Value = 24.05

$ws.Range("B4").Value = 110
$ws.Range("C4").Value = 2.6

$ws.Range("B5").Value = 110
$ws.Range("C5").Value = 2.6
